$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated values scraped from coinranking.com to the crypto table.
# Cells whose new value is a plain decimal number written as text (e.g.
# "315.54", "0.08195", "0.6370") are first switched to Text format so
# Excel keeps them as literal strings (preserving exact digits/trailing
# zeros) instead of silently parsing them into floating point numbers.

$ws.Range("D2").Value = '28.532.27'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.825.19'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.54'
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5107'
$ws.Range("E7").Value = '  -5.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3945'
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08195'
$ws.Range("E9").Value = '  +6.01%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.84'
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.356'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.13'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.562'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '1.820.67'
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001122'
$ws.Range("E17").Value = '  +3.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.95'
$ws.Range("E18").Value = '  +3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06659'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.105'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '28.569.95'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.41'
$ws.Range("E24").Value = '  +1.85%  '
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.05'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").Value = '2.029.76'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.408'
$ws.Range("E29").Value = '  -1.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.60'
$ws.Range("E30").Value = '  +1.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.115'
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  -3.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.772'
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.667'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07068'
$ws.Range("E35").Value = '  -4.00%  '
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02356'
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.254'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.768'
$ws.Range("E39").Value = '  -1.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6370'
$ws.Range("E40").Value = '  +1.23%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.182'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.400'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.66'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5977'
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("E46").Value = '  +0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.29'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.992'
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.196'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06946'
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.084'
$ws.Range("E51").Value = '  +4.07%  '
